$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.396.99"
$ws.Range("E2").Value = "  +4.29%  "

# Row 3
$ws.Range("D3").Value = "1.610.90"
$ws.Range("E3").Value = "  +2.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.54%  "

# Row 5
$ws.Range("D5").Value = "'213.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.61%  "

# Row 6
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("E7").Value = "  +2.35%  "

# Row 8
$ws.Range("D8").Value = "'0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.04%  "

# Row 9
$ws.Range("D9").Value = "'0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.27%  "

# Row 10
$ws.Range("D10").Value = "'18.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "

# Row 11
$ws.Range("D11").Value = "'0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.01%  "

# Row 12
$ws.Range("D12").Value = "1.837.14"
$ws.Range("E12").Value = "  +2.70%  "

# Row 13
$ws.Range("D13").Value = "1.610.63"
$ws.Range("E13").Value = "  +2.00%  "

# Row 14
$ws.Range("D14").Value = "'4.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
$ws.Range("D15").Value = "'0.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.59%  "

# Row 16
$ws.Range("D16").Value = "26.351.68"
$ws.Range("E16").Value = "  +4.08%  "

# Row 17
$ws.Range("D17").Value = "'61.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.29%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +2.78%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'208.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.12%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21
$ws.Range("D21").Value = "'4.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "

# Row 22
$ws.Range("D22").Value = "'9.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "

# Row 23
$ws.Range("D23").Value = "'6.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "

# Row 24
$ws.Range("D24").Value = "'1.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.22%  "

# Row 25
$ws.Range("D25").Value = "'142.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26
$ws.Range("E26").Value = "  -0.54%  "

# Row 27
$ws.Range("D27").Value = "'0.125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.38%  "

# Row 28
$ws.Range("D28").Value = "'15.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "

# Row 29
$ws.Range("D29").Value = "'6.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.40%  "

# Row 30
$ws.Range("E30").Value = "  +1.35%  "

# Row 31
$ws.Range("D31").Value = "'0.0474"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.48%  "

# Row 32
$ws.Range("D32").Value = "'3.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.28%  "

# Row 33
$ws.Range("D33").Value = "'3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "

# Row 34
$ws.Range("E34").Value = "  +1.82%  "

# Row 35
$ws.Range("E35").Value = "  +2.50%  "

# Row 36
$ws.Range("E36").Value = "  +8.33%  "

# Row 37
$ws.Range("D37").Value = "1.109.16"
$ws.Range("E37").Value = "  +1.70%  "

# Row 38
$ws.Range("E38").Value = "  -0.16%  "

# Row 39
$ws.Range("E39").Value = "  +0.76%  "

# Row 40
$ws.Range("D40").Value = "'0.789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.83%  "

# Row 41
$ws.Range("D41").Value = "'0.500"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "

# Row 42
$ws.Range("D42").Value = "'0.780"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.65%  "

# Row 43
$ws.Range("D43").Value = "1.746.21"
$ws.Range("E43").Value = "  +2.55%  "

# Row 44
$ws.Range("D44").Value = "'93.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

# Row 45
$ws.Range("E45").Value = "  +1.00%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.18%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0106"
$ws.Range("E47").Value = "  -4.98%  "

# Row 48
$ws.Range("D48").Value = "'53.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.84%  "

# Row 49
$ws.Range("D49").Value = "'0.0507"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "

# Row 50
$ws.Range("E50").Value = "  +0.37%  "
